$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = 1
$ws.Range("F25").Value = -2
$ws.Range("F31").Value = 1
$ws.Range("F34").Value = -1
$ws.Range("F35").Value = -1
$ws.Range("F41").Value = -1
$ws.Range("F43").Value = 1
$ws.Range("F48").Value = 5
$ws.Range("F49").Value = 5
$ws.Range("F51").Value = -4
$ws.Range("F54").Value = -4
$ws.Range("F55").Value = -4
$ws.Range("F56").Value = 1
$ws.Range("F57").Value = 1
$ws.Range("F58").Value = -1
$ws.Range("F61").Value = 0
$ws.Range("F63").Value = 1
$ws.Range("F69").Value = -1
$ws.Range("F73").Value = 0
$ws.Range("F76").Value = 1
$ws.Range("F80").Value = 0
$ws.Range("F81").Value = -1
